# Update the "Right Answer" values for question 11 (row 17) and
# question 15 (row 21) on Sheet1 from "a" to "b".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B17").Value = "b"
$ws.Range("B21").Value = "b"

# Move the active selection to B13, matching the saved cursor position.
$ws.Range("B13").Select()
